$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text values (e.g. "1.007", "27.913.88").
# Force text format on the whole Price column first so Excel does not
# silently coerce these strings into numeric values, then restore the
# default "Normal" style so no stray number format sticks to the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.913.88'
$ws.Range('D3').Value = '1.741.03'
$ws.Range('D4').Value = '1.007'
$ws.Range('D5').Value = '335.48'
$ws.Range('D7').Value = '0.3759'
$ws.Range('D8').Value = '0.3348'
$ws.Range('D9').Value = '45.17'
$ws.Range('D10').Value = '1.113'
$ws.Range('D11').Value = '0.07195'
$ws.Range('D12').Value = '1.004'
$ws.Range('D13').Value = '22.37'
$ws.Range('D14').Value = '6.157'
$ws.Range('D15').Value = '7.107'
$ws.Range('D16').Value = '1.746.08'
$ws.Range('D17').Value = '0.00001055'
$ws.Range('D18').Value = '0.06564'
$ws.Range('D19').Value = '79.03'
$ws.Range('D20').Value = '1.001'
$ws.Range('D21').Value = '16.81'
$ws.Range('D22').Value = '6.231'
$ws.Range('D23').Value = '27.931.40'
$ws.Range('D25').Value = '2.392'
$ws.Range('D26').Value = '154.00'
$ws.Range('D27').Value = '19.79'
$ws.Range('D28').Value = '2.315'
$ws.Range('D29').Value = '1.946.06'
$ws.Range('D30').Value = '131.40'
$ws.Range('D31').Value = '1.247'
$ws.Range('D32').Value = '4.028'
$ws.Range('D33').Value = '5.759'
$ws.Range('D34').Value = '0.08724'
$ws.Range('D36').Value = '0.6675'
$ws.Range('D37').Value = '0.02312'
$ws.Range('D38').Value = '0.06210'
$ws.Range('D39').Value = '5.146'
$ws.Range('D41').Value = '1.210'
$ws.Range('D42').Value = '1.439'
$ws.Range('D43').Value = '1.002'
$ws.Range('D44').Value = '7.941'
$ws.Range('D45').Value = '13.69'
$ws.Range('D46').Value = '3.821'
$ws.Range('D47').Value = '0.6031'
$ws.Range('D48').Value = '127.04'
$ws.Range('D49').Value = '2.015'
$ws.Range('D50').Value = '48.49'
$ws.Range('D51').Value = '0.07109'

$priceRange.Style = "Normal"

# Remaining text columns (Coin name, Link, Volume%) are not numeric-looking
# so they can be assigned directly without any coercion risk.
$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('E3').Value = '  -3.72%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  -4.06%  '
$ws.Range('E8').Value = '  -4.36%  '
$ws.Range('E9').Value = '  -6.98%  '
$ws.Range('E10').Value = '  -5.86%  '
$ws.Range('E11').Value = '  -4.84%  '
$ws.Range('E12').Value = '  +0.32%  '
$ws.Range('E13').Value = '  +0.97%  '
$ws.Range('E14').Value = '  -5.86%  '
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('E17').Value = '  -4.77%  '
$ws.Range('E18').Value = '  -2.28%  '
$ws.Range('E19').Value = '  -7.40%  '
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('E21').Value = '  -5.43%  '
$ws.Range('E22').Value = '  -5.25%  '
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('E24').Value = '  -6.85%  '
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('E26').Value = '  -0.63%  '
$ws.Range('E27').Value = '  -7.75%  '
$ws.Range('E28').Value = '  -8.48%  '
$ws.Range('E29').Value = '  -3.72%  '
$ws.Range('E30').Value = '  -3.89%  '
$ws.Range('E31').Value = '  -15.68%  '
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('E33').Value = '  -10.34%  '
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('E35').Value = '  -7.98%  '
$ws.Range('E36').Value = '  -3.73%  '
$ws.Range('E37').Value = '  -6.34%  '
$ws.Range('E38').Value = '  -5.32%  '
$ws.Range('E39').Value = '  -6.42%  '
$ws.Range('E40').Value = '  -5.61%  '
$ws.Range('E41').Value = '  -4.43%  '
$ws.Range('E42').Value = '  -11.21%  '
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E44').Value = '  -7.19%  '
$ws.Range('E45').Value = '  -6.45%  '
$ws.Range('E46').Value = '  -1.37%  '
$ws.Range('E47').Value = '  -6.37%  '
$ws.Range('E48').Value = '  -3.77%  '
$ws.Range('E49').Value = '  -7.06%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('E50').Value = '  +9.27%  '
$ws.Range('E51').Value = '  -1.95%  '
